$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Range("B22").Value = "SingleUseId25"
$ws.Range("C22").Value = "Default"
$ws.Range("D22").Value = "Left"
$ws.Range("E22").Value = "LTR"
$ws.Range("F22").Value = "Buzz"
